# "atualizei dados bibi e add"
#
# 1) The store labels in rows 4 and 5 were swapped: the row that used to
#    say "Bibi Cell Manauara" now says "Bibi Cell Vieiralves" and vice
#    versa (the daily figures that had been recorded under each name
#    stayed on their original row, so B4/C4 <-> B5/C5 swap too).
# 2) New revenue figures were added for several days (column D, and in
#    row 4 also columns B/C) for stores "Bibi Cell Mundi" (row 2),
#    "Bibi Cell Ponta Negra" (row 3), the two renamed stores (rows 4-5)
#    and the "total" row (row 6).
# 3) Column AG (the "total" column) is simply the sum of each row and is
#    updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the store names on rows 4 and 5 ---------------------------------
$ws.Range("A4").Value = "Bibi Cell Vieiralves"
$ws.Range("A5").Value = "Bibi Cell Manauara"

# --- Row 2: Bibi Cell Mundi -------------------------------------------------
$ws.Range("D2").Value = 17353.24
$ws.Range("AG2").Value = 34498.04

# --- Row 3: Bibi Cell Ponta Negra ------------------------------------------
$ws.Range("D3").Value = 1748.51
$ws.Range("AG3").Value = 8218.52

# --- Row 4: now Bibi Cell Vieiralves ---------------------------------------
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 4464
$ws.Range("D4").Value = 3201
$ws.Range("AG4").Value = 7665

# --- Row 5: now Bibi Cell Manauara ------------------------------------------
$ws.Range("B5").Value = 3340
$ws.Range("C5").Value = 1374
$ws.Range("D5").Value = 2934
$ws.Range("AG5").Value = 7648

# --- Row 6: total ------------------------------------------------------------
$ws.Range("D6").Value = 25236.75
$ws.Range("AG6").Value = 58029.56
